$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (Changed) date column (C) for rows 2-8
# from serial date 45221 (2023-10-22) to 45224 (2023-10-25).
foreach ($row in 2..8) {
    $ws.Cells.Item($row, 3).Value = 45224
}
